$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 24150
$ws.Range("D2").Value = 0.0168
$ws.Range("I2").Value = 4.14
$ws.Range("J2").Value = 84
$ws.Range("K2").Value = 84

$ws.Range("C3").Value = 101900
$ws.Range("D3").Value = 0.001
$ws.Range("D3").NumberFormat = "0.00%"
$ws.Range("I3").Value = 6.38

$ws.Range("C4").Value = 433500
$ws.Range("D4").Value = -0.0214
$ws.Range("I4").Value = 4.38
$ws.Range("J4").Value = 74
$ws.Range("K4").Value = 74

$ws.Range("C5").Value = 30200
$ws.Range("D5").Value = -0.0049
$ws.Range("I5").Value = 6.62
$ws.Range("J5").Value = 44
$ws.Range("K5").Value = 44

$ws.Range("C6").Value = 29250
$ws.Range("D6").Value = -0.0135
$ws.Range("I6").Value = 4.1
$ws.Range("J6").Value = 74
$ws.Range("K6").Value = 74

$ws.Range("C7").Value = 25200
$ws.Range("D7").Value = 0.004
$ws.Range("I7").Value = 4.76
$ws.Range("J7").Value = 69
$ws.Range("K7").Value = 69

$ws.Range("C8").Value = 10880
$ws.Range("D8").Value = 0.0112
$ws.Range("I8").Value = 4.73
$ws.Range("J8").Value = 89
$ws.Range("K8").Value = 89

$ws.Range("C9").Value = 88000
$ws.Range("D9").Value = 0.0057
$ws.Range("I9").Value = 3.41
$ws.Range("J9").Value = 82
$ws.Range("K9").Value = 82

$ws.Range("D10").Value = 0
$ws.Range("D10").NumberFormat = "0%"

$ws.Range("C11").Value = 127600
$ws.Range("D11").Value = -0.0139
$ws.Range("I11").Value = 5.33
$ws.Range("J11").Value = 81
$ws.Range("K11").Value = 81

$ws.Range("C12").Value = 19500
$ws.Range("D12").Value = -0.0051
$ws.Range("I12").Value = 4.87

$ws.Range("C13").Value = 70000
$ws.Range("D13").Value = -0.0085
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 79
$ws.Range("K13").Value = 79

$ws.Range("C14").Value = 56100
$ws.Range("D14").Value = 0.0108
$ws.Range("I14").Value = 6.31
$ws.Range("J14").Value = 73
$ws.Range("K14").Value = 73

$ws.Range("C15").Value = 85500
$ws.Range("D15").Value = 0.0339
$ws.Range("I15").Value = 6.43
$ws.Range("J15").Value = 91
$ws.Range("K15").Value = 91

$ws.Range("C16").Value = 19400
$ws.Range("D16").Value = 0.0015
$ws.Range("I16").Value = 5.49

$ws.Range("C17").Value = 50700
$ws.Range("D17").Value = 0.006
$ws.Range("D17").NumberFormat = "0.00%"
$ws.Range("I17").Value = 5.52
$ws.Range("J17").Value = 72
$ws.Range("K17").Value = 72

$ws.Range("C18").Value = 20100
$ws.Range("D18").Value = 0.0025
$ws.Range("I18").Value = 6.12

$ws.Range("C19").Value = 53700
$ws.Range("D19").Value = -0.011
$ws.Range("I19").Value = 3.72
$ws.Range("J19").Value = 85
$ws.Range("K19").Value = 85

$ws.Range("C20").Value = 14590
$ws.Range("D20").Value = 0.0021
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = 76

$ws.Range("C21").Value = 147200
$ws.Range("D21").Value = 0.0628
$ws.Range("I21").Value = 3.67
$ws.Range("J21").Value = 97
$ws.Range("K21").Value = 97

$ws.Range("C22").Value = 41300
$ws.Range("D22").Value = -0.0236
$ws.Range("I22").Value = 3.52
$ws.Range("J22").Value = 44
$ws.Range("K22").Value = 44

$ws.Range("D23").Value = 0
$ws.Range("D23").NumberFormat = "0%"

$ws.Range("C24").Value = 48050
$ws.Range("I24").Value = 5.62
$ws.Range("J24").Value = 65
$ws.Range("K24").Value = 65

$ws.Range("C25").Value = 84900
$ws.Range("D25").Value = -0.0105
$ws.Range("I25").Value = 4.24
$ws.Range("J25").Value = 83
$ws.Range("K25").Value = 83

$ws.Range("C26").Value = 113500
$ws.Range("D26").Value = -0.013
$ws.Range("I26").Value = 2.8
$ws.Range("J26").Value = 85
$ws.Range("K26").Value = 85

$ws.Range("C27").Value = 14660
$ws.Range("D27").Value = 0.0027
$ws.Range("I27").Value = 4.43

$ws.Range("C28").Value = 13910
$ws.Range("D28").Value = -0.01
$ws.Range("I28").Value = 3.59
$ws.Range("J28").Value = 84
$ws.Range("K28").Value = 84

$ws.Range("C29").Value = 22050
$ws.Range("D29").Value = -0.02
$ws.Range("I29").Value = 4.51
$ws.Range("J29").Value = 81
$ws.Range("K29").Value = 81

$ws.Range("C30").Value = 25200
$ws.Range("D30").Value = -0.0059
$ws.Range("I30").Value = 4.76
$ws.Range("J30").Value = 90
$ws.Range("K30").Value = 90

$ws.Range("G18").Select()
